{"js": "// The author's edit inserts the word \"El \" right before\n// \"Se\u00f1or Notario, esta\" in the closing paragraph of the act\n// (\"...compareciente por mi Se\u00f1or Notario, esta...\" ->\n//  \"...compareciente por mi El Se\u00f1or Notario, esta...\").\n//\n// (Word's own \"_GoBack\" bookmark - which tracks the location of the\n// last edit - also shifts from the start of that sentence to right\n// before \"Se\u00f1or Notario, esta\" as a side effect of typing there; we\n// replicate that too.)\n\nconst body = context.document.body;\n\n// Locate the unique target phrase and insert \"El \" right before it.\nconst results = body.search(\"Se\u00f1or Notario, esta\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Target phrase \"Se\u00f1or Notario, esta\" not found.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\"El \", Word.InsertLocation.before);\nawait context.sync();\n\n// Move the \"_GoBack\" bookmark so it again marks the last-edited spot,\n// i.e. right before \"Se\u00f1or Notario, esta\".\nconst goBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\ngoBack.load(\"isNullObject\");\nawait context.sync();\n\nif (!goBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\nconst results2 = body.search(\"Se\u00f1or Notario, esta\", { matchCase: true, matchWholeWord: false });\nresults2.load(\"items\");\nawait context.sync();\n\nconst newAnchor = results2.items[0].getRange(Word.RangeLocation.start);\nnewAnchor.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author's edit inserts the word \"El \" right before\n# \"Se\u00f1or Notario, esta\" in the closing paragraph of the act\n# (\"...compareciente por mi Se\u00f1or Notario, esta...\" ->\n#  \"...compareciente por mi El Se\u00f1or Notario, esta...\").\n#\n# (Word's own \"_GoBack\" bookmark - which tracks the location of the\n# last edit - also shifts from the start of that sentence to right\n# before \"Se\u00f1or Notario, esta\" as a side effect of typing there; we\n# replicate that too.)\n\n$d = $word.ActiveDocument\n\n# Locate the unique target phrase and insert \"El \" right before it.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Se\u00f1or Notario, esta\")\nif (-not $found) {\n    throw 'Target phrase \"Se\u00f1or Notario, esta\" not found.'\n}\n$rng.Collapse(1)\n$rng.InsertBefore(\"El \")\n\n# Move the \"_GoBack\" bookmark so it again marks the last-edited spot,\n# i.e. right before \"Se\u00f1or Notario, esta\".\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\"Se\u00f1or Notario, esta\")\nif (-not $found2) {\n    throw 'Target phrase \"Se\u00f1or Notario, esta\" not found (pass 2).'\n}\n$rng2.Collapse(1)\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$d.Bookmarks.Add(\"_GoBack\", $rng2)\n"}
